$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Worksheet view state (best effort, cosmetic) -----------------------
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A19").Select()
$win.FreezePanes = $true
$win.Zoom = 144

# --- Row 12 header (player names for "Contest 5") ------------------------
$ws.Range("S12").Value = "Sampath"
$ws.Range("V12").Value = "Jayanth"

# --- Row 14: PBKS vs DC ---------------------------------------------------
$ws.Range("C14").Value = "PBKS vs DC"
$ws.Range("E14").Value = 60
$ws.Range("H14").Value = 30
$ws.Range("K14").Value = 70
$ws.Range("N14").Value = 20
$ws.Range("Q14").Value = 40
$ws.Range("T14").Value = 50
$ws.Range("W14").Value = 100
$ws.Range("Z14").Value = 20
# Tied ranks (N14 = Z14 = 20) -> manually averaged points, formula cleared
$ws.Range("M14").Value = -22.5
$ws.Range("Y14").Value = -22.5

# --- Row 15: KKR vs SRH ----------------------------------------------------
$ws.Range("C15").Value = "KKR vs SRH"
$ws.Range("E15").Value = 40
$ws.Range("H15").Value = 50
$ws.Range("K15").Value = 60
$ws.Range("N15").Value = 100
$ws.Range("Q15").Value = 70
$ws.Range("T15").Value = 20
$ws.Range("W15").Value = 30
$ws.Range("Z15").Value = 20
# Tied ranks (T15 = Z15 = 20) -> manually averaged points, formula cleared
$ws.Range("S15").Value = -22.5
$ws.Range("Y15").Value = -22.5

# --- Row 16: RR vs LSG -------------------------------------------------
$ws.Range("C16").Value = "RR vs LSG"
$ws.Range("E16").Value = 0
$ws.Range("H16").Value = 50
$ws.Range("K16").Value = 60
$ws.Range("N16").Value = 30
$ws.Range("Q16").Value = 100
$ws.Range("T16").Value = 40
$ws.Range("W16").Value = 70
$ws.Range("Z16").Value = 20

# --- Row 17: GT vs MI --------------------------------------------------
$ws.Range("C17").Value = "GT vs MI"
$ws.Range("E17").Value = 50
$ws.Range("H17").Value = 20
$ws.Range("K17").Value = 70
$ws.Range("N17").Value = 100
$ws.Range("Q17").Value = 30
$ws.Range("T17").Value = 0
$ws.Range("W17").Value = 40
$ws.Range("Z17").Value = 60

# --- Row 18: next match label only -------------------------------------
$ws.Range("C18").Value = "RCB vs PBKS"

# --- Bottom summary row (74): point totals now pull live from row 12 ----
$ws.Range("T74").Formula = "=S12"
$ws.Range("W74").Formula = "=V12"

# --- Restore final selection/scroll position ----------------------------
$ws.Range("A68").Select()
$ws.Range("AA75").Select()
